# Update odds values in the weekly FlashScore games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G4"   = 1.95
    "H4"   = 2.9
    "M4"   = 1.14
    "N4"   = 5.5
    "O4"   = 1.67
    "P4"   = 2.1
    "Q4"   = 3.1
    "R4"   = 1.36
    "S4"   = 1.73
    "T4"   = 2
    "X4"   = 7
    "AD4"  = 6.5
    "AH4"  = 23
    "AO4"  = 12
    "AU4"  = 12

    "G5"   = 2.25
    "H5"   = 3.1
    "J5"   = 3.1
    "K5"   = 1.91
    "M5"   = 1.13
    "N5"   = 6
    "O5"   = 1.53
    "P5"   = 2.38
    "U5"   = 2.2
    "V5"   = 1.62
    "AS5"  = 301
    "AW5"  = 5.5
    "AX5"  = 23

    "M9"   = 1.07
    "N9"   = 9

    "G11"  = 1.95
    "I11"  = 4.2
    "L11"  = 4.5
    "X11"  = 8.5
    "AA11" = 17
    "AH11" = 21
    "AO11" = 11
    "AX11" = 23
    "BB11" = 301
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
